$d = $word.ActiveDocument

# The "_GoBack" bookmark currently sits at the end of the last journal
# entry ("Created an error message..."). It needs to move to the end of
# the new last entry we are about to append, so drop it here and re-add
# it in the new final paragraph below.
$oldGoBack = $d.Bookmarks.Item("_GoBack")
$oldGoBack.Delete()

# Find the paragraph that currently ends the document's journal content:
# "Created an error message for when the program attempts to divide by 0"
$needle = "Created an error message for when the program attempts to divide by 0"
$lastEntry = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*$needle*") {
        $lastEntry = $candidate
    }
}

$anchor = $lastEntry.Range

# Add two new paragraphs after it: one for the date, one for the new
# journal entry.
$anchor.InsertParagraphAfter()
$anchor.InsertParagraphAfter()

$dateParaIndex = $lastEntry.Index + 1
$entryParaIndex = $lastEntry.Index + 2

$dateParagraph = $d.Paragraphs.Item($dateParaIndex)
$dateParagraph.Range.InsertAfter("12/06/19")

$entryParagraph = $d.Paragraphs.Item($entryParaIndex)
$entryXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Set the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>combobox</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> for the hours per week to editable</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$entryParagraph.Range.InsertXML($entryXml)
